$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 270
$ws.Range("F4").Value = 280
$ws.Range("F5").Value = 2888
$ws.Range("F8").Value = 2244
$ws.Range("F9").Value = 1441
$ws.Range("F10").Value = 1441
$ws.Range("F11").Value = 39
$ws.Range("F12").Value = 450
$ws.Range("F14").Value = 2581
$ws.Range("F16").Value = 1399
$ws.Range("F17").Value = 4893
$ws.Range("F19").Value = 5358
$ws.Range("F20").Value = 5358
$ws.Range("F21").Value = 1921
$ws.Range("F22").Value = 2936
$ws.Range("F23").Value = 3354
$ws.Range("F24").Value = 192
$ws.Range("F25").Value = 1626
$ws.Range("F26").Value = 271
$ws.Range("F28").Value = 137
$ws.Range("F29").Value = 4
$ws.Range("F30").Value = 327
$ws.Range("F31").Value = 1045
$ws.Range("F32").Value = 2136
$ws.Range("F33").Value = 3
$ws.Range("F34").Value = 129
$ws.Range("F35").Value = 304
$ws.Range("F36").Value = 788
$ws.Range("F37").Value = 166
$ws.Range("F38").Value = 372
$ws.Range("F39").Value = 445

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 105
$ws.Range("F6").Value = 59
$ws.Range("F12").Value = 22
$ws.Range("F13").Value = 2
$ws.Range("F15").Value = 15
$ws.Range("F16").Value = 45

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 105
$ws.Range("F6").Value = 59
$ws.Range("F7").Value = 270
$ws.Range("F9").Value = 280
$ws.Range("F10").Value = 2888
$ws.Range("F12").Value = 2244
$ws.Range("F13").Value = 1441
$ws.Range("F14").Value = 1441
$ws.Range("F16").Value = 39
$ws.Range("F17").Value = 450
$ws.Range("F20").Value = 2581
$ws.Range("F21").Value = 1399
$ws.Range("F23").Value = 22
$ws.Range("F24").Value = 2
$ws.Range("F26").Value = 4893
$ws.Range("F28").Value = 5358
$ws.Range("F29").Value = 5358
$ws.Range("F30").Value = 1921
$ws.Range("F31").Value = 2936
$ws.Range("F32").Value = 3354
$ws.Range("F33").Value = 192
$ws.Range("F34").Value = 45
$ws.Range("F36").Value = 1626
$ws.Range("F38").Value = 271
$ws.Range("F40").Value = 137
$ws.Range("F41").Value = 4
$ws.Range("F42").Value = 327
$ws.Range("F44").Value = 2136
$ws.Range("F45").Value = 3
$ws.Range("F46").Value = 129
$ws.Range("F47").Value = 304
$ws.Range("F48").Value = 788
$ws.Range("F49").Value = 166
$ws.Range("F50").Value = 372
$ws.Range("F51").Value = 445
